# update form filling to use reflection
# Fills in the "Personal Info Template" sheet with sample personal data,
# and leaves the active sheet/selection pointed at the sheet & cell the
# author was last working on.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Personal Info Template")

# --- Personal info fields ---
$ws1.Range("B2").Value  = "Hu"            # Last name
$ws1.Range("B4").Value  = "Jack"          # First name
$ws1.Range("B5").Value  = 123456789       # Identifying number
$ws1.Range("B6").Value  = "220 Patina Green SW"  # Home address
$ws1.Range("B7").Value  = "Calgary"       # City
$ws1.Range("B8").Value  = "Canada"        # Country
$ws1.Range("B9").Value  = "Alberta"       # Province/State
$ws1.Range("B10").Value = "T3H3C7"        # Postal code
$ws1.Range("B11").Value = 3               # Filing status
$ws1.Range("B12").Value = "Allan Hu"      # Child's name(s)

# --- Spouse info fields ---
$ws1.Range("B14").Value = "test1"         # Spouse's last name
$ws1.Range("B16").Value = "test3"         # Spouse's first name
$ws1.Range("B15").Value = "t"             # Spouse's middle initial
$ws1.Range("B17").Value = 123456789       # Spouse's identifying number

# Bring "Personal Info Template" to the front and leave the selection on B15,
# matching where the author left off editing.
$ws1.Activate()
$ws1.Range("B15").Select() | Out-Null
